# fix save grap and train process
# Updates recomputed metric columns G (precision?), H, I on Sheet1
# for the rows that have non-trivial (non-all-zero) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.03610305060647555
$ws.Cells.Item(2, 8).Value = 0.02962802687438672
$ws.Cells.Item(2, 9).Value = 0.8966768717453046
$ws.Cells.Item(3, 7).Value = 0.04206075011517957
$ws.Cells.Item(3, 8).Value = 0.02941297523609051
$ws.Cells.Item(3, 9).Value = 0.7792759239402862
$ws.Cells.Item(4, 7).Value = 0.04367587408253915
$ws.Cells.Item(4, 8).Value = 0.03851780101009961
$ws.Cells.Item(4, 9).Value = 0.849636719234776
$ws.Cells.Item(5, 7).Value = 0.01279733275079968
$ws.Cells.Item(5, 8).Value = 0.009877648428077685
$ws.Cells.Item(5, 9).Value = 0.9841915360654097
$ws.Cells.Item(6, 7).Value = 0.02863196580284044
$ws.Cells.Item(6, 8).Value = 0.02217154059815238
$ws.Cells.Item(6, 9).Value = 0.93238304070893
$ws.Cells.Item(7, 7).Value = 0.00072164456540591
$ws.Cells.Item(7, 8).Value = 0.00072164456540591
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(8, 7).Value = 0.04854975944118319
$ws.Cells.Item(8, 8).Value = 0.03972664576835257
$ws.Cells.Item(8, 9).Value = 0.7401901702060596
$ws.Cells.Item(9, 7).Value = 0.03870507238061709
$ws.Cells.Item(9, 8).Value = 0.03227888378133689
$ws.Cells.Item(9, 9).Value = 0.9041125588155597
$ws.Cells.Item(10, 7).Value = 0.008935991330989443
$ws.Cells.Item(10, 8).Value = 0.006808702845077137
$ws.Cells.Item(10, 9).Value = 0.9922561561314475
$ws.Cells.Item(11, 7).Value = 0.02754970086832659
$ws.Cells.Item(11, 8).Value = 0.02231273774802066
$ws.Cells.Item(11, 9).Value = 0.9361703187624111
$ws.Cells.Item(12, 7).Value = 0.03131467268923646
$ws.Cells.Item(12, 8).Value = 0.0271425883500411
$ws.Cells.Item(12, 9).Value = 0.9200263497335257
$ws.Cells.Item(14, 7).Value = 0.1070204674541322
$ws.Cells.Item(14, 8).Value = 0.07998116939477266
$ws.Cells.Item(14, 9).Value = 0.6016821687551963
$ws.Cells.Item(15, 7).Value = 0.01253458574461442
$ws.Cells.Item(15, 8).Value = 0.01137535145395804
$ws.Cells.Item(15, 9).Value = 0.9812857651222114
$ws.Cells.Item(16, 7).Value = 0.04995075777609275
$ws.Cells.Item(16, 8).Value = 0.03651754629601307
$ws.Cells.Item(16, 9).Value = 0.8048082600575938
$ws.Cells.Item(17, 7).Value = 0.03728957186256869
$ws.Cells.Item(17, 8).Value = 0.03038515542027997
$ws.Cells.Item(17, 9).Value = 0.9025796501827765
$ws.Cells.Item(19, 7).Value = 0.009419103519648418
$ws.Cells.Item(19, 8).Value = 0.007161817792903626
$ws.Cells.Item(19, 9).Value = 0.9919993158397526
$ws.Cells.Item(20, 7).Value = 0.2260000603956003
$ws.Cells.Item(20, 8).Value = 0.1860703052033772
$ws.Cells.Item(20, 9).Value = -0.5428640119040629
$ws.Cells.Item(21, 7).Value = 0.01476642655127117
$ws.Cells.Item(21, 8).Value = 0.01265350669783064
$ws.Cells.Item(21, 9).Value = 0.9795193378140926
$ws.Cells.Item(23, 7).Value = 0.02991437724315698
$ws.Cells.Item(23, 8).Value = 0.0231693272788059
$ws.Cells.Item(23, 9).Value = 0.9356243603816317
$ws.Cells.Item(25, 7).Value = 0.0727321099105632
$ws.Cells.Item(25, 8).Value = 0.06249923890352882
$ws.Cells.Item(25, 9).Value = 0.7323485042637632
$ws.Cells.Item(26, 7).Value = 0.03997878757854137
$ws.Cells.Item(26, 8).Value = 0.03033079754323986
$ws.Cells.Item(26, 9).Value = 0.8892899343886684
$ws.Cells.Item(27, 7).Value = 0.01492286784809646
$ws.Cells.Item(27, 8).Value = 0.01216345917277525
$ws.Cells.Item(27, 9).Value = 0.9852278508375966
$ws.Cells.Item(28, 7).Value = 0.04159517268588053
$ws.Cells.Item(28, 8).Value = 0.03576272667993004
$ws.Cells.Item(28, 9).Value = 0.8799709221938732
$ws.Cells.Item(29, 7).Value = 0.03381573516544372
$ws.Cells.Item(29, 8).Value = 0.02836332156304972
$ws.Cells.Item(29, 9).Value = 0.9002339580999384
$ws.Cells.Item(30, 7).Value = 0.01382792199348841
$ws.Cells.Item(30, 8).Value = 0.01159478279518422
$ws.Cells.Item(30, 9).Value = 0.9800448503322999
$ws.Cells.Item(31, 7).Value = 0.02283248362889998
$ws.Cells.Item(31, 8).Value = 0.01549194593669462
$ws.Cells.Item(31, 9).Value = 0.9298935281585758
$ws.Cells.Item(32, 7).Value = 0.01696810558439541
$ws.Cells.Item(32, 8).Value = 0.01493147344804608
$ws.Cells.Item(32, 9).Value = 0.9691069322266278
$ws.Cells.Item(34, 7).Value = 0.09109114755556109
$ws.Cells.Item(34, 8).Value = 0.07217103055672187
$ws.Cells.Item(34, 9).Value = 0.6369094188915406
$ws.Cells.Item(38, 7).Value = 0.2014864056314926
$ws.Cells.Item(38, 8).Value = 0.1710316697266309
$ws.Cells.Item(38, 9).Value = -0.4230806214895253
$ws.Cells.Item(40, 7).Value = 0.002734970331026548
$ws.Cells.Item(40, 8).Value = 0.002368550044713347
$ws.Cells.Item(40, 9).Value = 0.9991629287485503
$ws.Cells.Item(41, 7).Value = 0.008786255820602645
$ws.Cells.Item(41, 8).Value = 0.007668715473302345
$ws.Cells.Item(41, 9).Value = 0.991792980826775
